$d = $word.ActiveDocument

# The document originally has 6 body paragraphs (+ sectPr):
#   1. "4-bytes-saga"
#   2. (empty)
#   3. "spel 1: Balanceren. ..."
#   4. "spel 2: Whack a mole. ..."
#   5. "spel 3: Jumping game on LCD screen"
#   6. "spel 4: Boss fight"
#
# Paragraphs 3-6 are each rewritten/expanded into 2 (or 3, for the
# original paragraph 3) new paragraphs. We use Range.InsertXML, which
# replaces the contents of the target Range with the supplied WordML,
# to swap each original paragraph for its replacement paragraph(s) in
# one shot (this also lets us control run splits and the <w:proofErr/>
# spell-check markers exactly). We go from the last paragraph back to
# the first so earlier paragraph indices stay valid while we work.

# --- Paragraph 6: "spel 4: Boss fight" -> "Spel 4: " / "ontsnap het doolhof..."
$p6 = $d.Paragraphs.Item(6)
$xml6 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Spel </w:t></w:r><w:r><w:t>4</w:t></w:r><w:r><w:t xml:space="preserve">: </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>ontsnap het doolhof, ren door het doolhof met dezelfde knoppen als de map</w:t></w:r><w:r><w:t>, gebruik de het nieuwe gereedschap (knop 5 van links voor de hamer en knop 6 voor de poolstok) om door obstakels te komen</w:t></w:r><w:r><w:t xml:space="preserve"> en ontsnap eindelijk!</w:t></w:r></w:p>
'@
$p6.Range.InsertXML($xml6)

# --- Paragraph 5: "spel 3: Jumping game on LCD screen" -> "spel 3: " / "Parcour!..."
$p5 = $d.Paragraphs.Item(5)
$xml5 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">spel 3: </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:r><w:t>Parcour</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>! Wissel tussen de bovenkant en onderkant van de rijen op het lcd scherm en bereik het einde om de poolstok te krijgen waarmee je over gaten kunt springen.</w:t></w:r></w:p>
'@
$p5.Range.InsertXML($xml5)

# --- Paragraph 4: "spel 2: Whack a mole. ..." -> "spel 2: " / "Whack a mole. ... hamer ..."
$p4 = $d.Paragraphs.Item(4)
$xml4 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">spel </w:t></w:r><w:r><w:t>2</w:t></w:r><w:r><w:t xml:space="preserve">: </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:r><w:t>Whack</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>mole</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>. Druk de verschillende lampjes zullen willekeurig gaan branden. Druk de knoppen onder de lampjes die beginnen te branden op tijd aan om punten te scoren. Elke fout haalt 2 punten weg. Bereik de 2</w:t></w:r><w:r><w:t>5</w:t></w:r><w:r><w:t xml:space="preserve"> punten en win. Het spel zal moeilijker worden hoe dichter je bij het einde komt</w:t></w:r><w:r><w:t xml:space="preserve">! </w:t></w:r><w:r><w:br/><w:t>Als je het spel voltooid krijg je de hamer die obstakels in het doolhof kan slopen.</w:t></w:r></w:p>
'@
$p4.Range.InsertXML($xml4)

# --- Paragraph 3: "spel 1: Balanceren. ..." -> intro / "spel 1: " / "Balanceren. ... compas ..."
$p3 = $d.Paragraphs.Item(3)
$xml3 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Navigeer rond de map om de spellen te vinden. Gebruik de eerste knop links of naar links te bewegen, de tweede knop om naar rechts te bewegen, de derde knop om omhoog te bewegen en de vierde knop om omlaag te bewegen. </w:t></w:r><w:r><w:t>Voltooi alle spellen om items te krijgen en via het doolhof te kunnen ontsnappen!</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">spel </w:t></w:r><w:r><w:t>1</w:t></w:r><w:r><w:t xml:space="preserve">: </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Balanceren. Gebruik de potentiometer om het brandende lampje in het midden te houden voor 10 seconden</w:t></w:r><w:r><w:t>.</w:t></w:r><w:r><w:br/><w:t xml:space="preserve">Als je het spel voltooid heb krijg je het </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>compas</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> dat je richting het einde</w:t></w:r><w:r><w:t xml:space="preserve"> van het doolhof</w:t></w:r><w:r><w:t xml:space="preserve"> wijst</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p>
'@
$p3.Range.InsertXML($xml3)

Write-Output ("Paragraph count: " + $d.Paragraphs.Count)
